# Insert a new weekly price record for "Papa" (Cardinal, 1a (cosecha)) dated
# 2021-12-21 (serial 44551) into the "Terminal La Palmera de La Serena - Papa"
# data sheet, right before the existing row for 2021-07-09 (old row 247).
# All rows from the old row 247 through the old last row (305) shift down by
# one row to make room, which is achieved with a native row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 247..305 down to 248..306, duplicating formatting of row 247.
$ws.Rows(247).Insert()

# Populate the newly inserted row 247 with the new record's data.
$ws.Range("A247").Value = 8
$ws.Range("B247").Value = "Terminal La Palmera de La Serena"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44551
$ws.Range("E247").Value = 4
$ws.Range("F247").Value = 100114001
$ws.Range("G247").Value = "Papa"
$ws.Range("H247").Value = "Cardinal"
$ws.Range("I247").Value = "1a (cosecha)"
$ws.Range("J247").Value = 2400
$ws.Range("K247").Value = 12000
$ws.Range("L247").Value = 12500
$ws.Range("M247").Value = 12250
$ws.Range("N247").Value = "`$/saco 25 kilos"
$ws.Range("O247").Value = "Provincia del Elquí"
$ws.Range("P247").Value = 490
$ws.Range("Q247").Value = 25
$ws.Range("R247").Value = "Hortaliza"
